# Add "clubs" class/sheet and "student_clubs" sheet, wire up headers and
# summary cells (num_clubs / num_students_club), mirroring the other
# "register / view details / add" helper sheets already in the workbook.

$wb = $excel.ActiveWorkbook

# --- add the "clubs" sheet right after "student_courses" -------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$clubs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$clubs.Name = "clubs"

$clubs.Range("A1").Value = "Club_ID"
$clubs.Range("B1").Value = "Club_name"
$clubs.Range("C1").Value = "Subject"
$clubs.Range("D1").Value = "Description"

$clubs.Range("F3").Value = "num_clubs"
$clubs.Range("F4").Value = 0

# --- add the "student_clubs" sheet right after "clubs" ---------------------
$studentClubs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $clubs)
$studentClubs.Name = "student_clubs"

$studentClubs.Range("A1").Value = "Student_ID"
$studentClubs.Range("B1").Value = "fname"
$studentClubs.Range("C1").Value = "lname"
$studentClubs.Range("D1").Value = "Club_ID"
$studentClubs.Range("E1").Value = "Club_name"

$studentClubs.Range("H3").Value = "num_students_club"
$studentClubs.Range("H4").Value = 0

# the "register" column (A) gets the same bold-ish style used elsewhere
$studentClubs.Range("A1:A5").Font.Bold = $true

# the new "clubs" sheet is the active tab, matching the workbook's tracked
# selection after the edit
$clubs.Activate()
